$wb = $excel.ActiveWorkbook

# --- Update the Date value on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-09-12T13:34:32+00:00"

# --- Update the Definition column on the Concepts sheet with properly-cased text ---
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("C3").Value = "Missing - Restricted Access"
$concepts.Range("C4").Value = "Missing - Not Provided"
$concepts.Range("C5").Value = "Missing - Not Collected"
